# "매출액 업로드 첨부파일.xlsx" - update the header row on Sheet1 to match the
# new SaleManager domain fields (commit: "Feat:Add SaleManager domain and
# Create fileUpload API").
#
# Old headers: 회사명 | 사업자번호 | 매출일(년도) | 매출일(월) | 매출액
# New headers: 회사명 | 사업자번호 | 매출액       | 영업이익   | 당기순이익 | 기준일자(년월)
#
# Columns C/D drop the split year/month sales-date fields; the old "매출액"
# column slides left to C, two new metric columns (영업이익, 당기순이익) are
# added, and a single combined "기준일자(년월)" date column is appended as F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "매출액"
$ws.Range("D1").Value = "영업이익"
$ws.Range("E1").Value = "당기순이익"
$ws.Range("F1").Value = "기준일자(년월)"

# The edited workbook was left with the newly added F1 header cell selected.
$ws.Range("F1").Select()
